$wb = $excel.ActiveWorkbook

# --- Sheet "Detalle": individual status corrections ---
$detalle = $wb.Worksheets.Item("Detalle")

$detalle.Cells.Item(202, 2).Value = 4
$detalle.Cells.Item(202, 3).Value = "Migrado"
$detalle.Cells.Item(229, 2).Value = 7
$detalle.Cells.Item(229, 3).Value = "Por Reprogramar"
$detalle.Cells.Item(589, 2).Value = 7
$detalle.Cells.Item(589, 3).Value = "Por Reprogramar"
$detalle.Cells.Item(595, 2).Value = 4
$detalle.Cells.Item(595, 3).Value = "Migrado"
$detalle.Cells.Item(874, 2).Value = 4
$detalle.Cells.Item(874, 3).Value = "Migrado"

# --- Sheet "Detalle": rows 2002-2050 reordered/updated ---
$rows2002to2050 = @(
    @(2002, 79166, 9, "En Bodega", 45894),
    @(2003, 79162, 9, "En Bodega", 45890),
    @(2004, 80084, 2, "Alistamiento", 45936),
    @(2005, 79332, 2, "Alistamiento", 45934),
    @(2006, 79334, 2, "Alistamiento", 45910),
    @(2007, 79446, 9, "En Bodega", 45881),
    @(2008, 79440, 2, "Alistamiento", 45901),
    @(2009, 79588, 9, "En Bodega", 45884),
    @(2010, 79671, 2, "Alistamiento", 45919),
    @(2011, 79674, 2, "Alistamiento", 45919),
    @(2012, 79809, 9, "En Bodega", 45880),
    @(2013, 79966, 9, "En Bodega", 45895),
    @(2014, 79164, 2, "Alistamiento", 45901),
    @(2015, 79163, 2, "Alistamiento", 45923),
    @(2016, 80085, 2, "Alistamiento", 45912),
    @(2017, 79432, 2, "Alistamiento", 45919),
    @(2018, 79336, 2, "Alistamiento", 45910),
    @(2019, 79436, 9, "En Bodega", 45888),
    @(2020, 79441, 2, "Alistamiento", 45936),
    @(2021, 79437, 9, "En Bodega", 45895),
    @(2022, 79589, 9, "En Bodega", 45890),
    @(2023, 79672, 2, "Alistamiento", 45959),
    @(2024, 79675, 2, "Alistamiento", 45929),
    @(2025, 79962, 9, "En Bodega", 45889),
    @(2026, 80082, 2, "Alistamiento", 45922),
    @(2027, 12543, 2, "Alistamiento", 45931),
    @(2028, 79165, 2, "Alistamiento", 45912),
    @(2029, 79204, 2, "Alistamiento", 45923),
    @(2030, 80318, 9, "En Bodega", 45881),
    @(2031, 79333, 9, "En Bodega", 45894),
    @(2032, 79435, 9, "En Bodega", 45890),
    @(2033, 79439, 2, "Alistamiento", 45929),
    @(2034, 79587, 9, "En Bodega", 45894),
    @(2035, 79669, 2, "Alistamiento", 45959),
    @(2036, 79673, 2, "Alistamiento", 45909),
    @(2037, 79677, 2, "Alistamiento", 45909),
    @(2038, 79808, 2, "Alistamiento", 45909),
    @(2039, 79965, 9, "En Bodega", 45896),
    @(2040, 79167, 2, "Alistamiento", 45904),
    @(2041, 79161, 2, "Alistamiento", 45933),
    @(2042, 80086, 9, "En Bodega", 45898),
    @(2043, 79433, 2, "Alistamiento", 45960),
    @(2044, 79434, 2, "Alistamiento", 45911),
    @(2045, 79438, 2, "Alistamiento", 45908),
    @(2046, 79442, 2, "Alistamiento", 45898),
    @(2047, 79586, 2, "Alistamiento", 45898),
    @(2048, 79668, 9, "En Bodega", 45884),
    @(2049, 79590, 2, "Alistamiento", 45898),
    @(2050, 79676, 2, "Alistamiento", 45950)
)

foreach ($row in $rows2002to2050) {
    $r = $row[0]
    $detalle.Cells.Item($r, 1).Value = $row[1]
    $detalle.Cells.Item($r, 2).Value = $row[2]
    $detalle.Cells.Item($r, 3).Value = $row[3]
    $detalle.Cells.Item($r, 4).Value = $row[4]
}

# --- Sheet "Resumen_por_estado": updated summary counts ---
$resumen = $wb.Worksheets.Item("Resumen_por_estado")
$resumen.Cells.Item(3, 2).Value = 305
$resumen.Cells.Item(4, 2).Value = 10
$resumen.Cells.Item(5, 1).Value = "Por Reprogramar"
$resumen.Cells.Item(5, 2).Value = 2
